$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Agua quantity (B2) becomes a text "20" instead of a number
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "20"

# Fanta quantity (B5) changes from "10" to "0"
$ws.Range("B5").Value = "0"

# Vino quantity (B6) changes from "30" to "25"
$ws.Range("B6").Value = "25"

# Tang quantity (B7) changes from "5" to "0"
$ws.Range("B7").Value = "0"

# New row 8: Pitusas
$ws.Range("A8").Value = "Pitusas"
$ws.Range("B8").Value = "20"
$ws.Range("C8").Value = "15"
